$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 671.375
$ws.Range("I9").Value = 592
$ws.Range("J9").Value = 733.1111
$ws.Range("K9").Value = 592
$ws.Range("L9").Value = 733.1111
$ws.Range("M9").Value = -423
$ws.Range("N9").Value = -1071.1111

$ws.Range("H18").Value = 297.25
$ws.Range("J18").Value = 242.5
$ws.Range("L18").Value = 242.5
$ws.Range("N18").Value = -810.5

$ws.Range("H19").Value = 2415.2222
$ws.Range("I19").Value = 2208.5454
$ws.Range("K19").Value = 2208.5454
$ws.Range("M19").Value = -2033.5454

$ws.Range("H28").Value = 1125.5385
$ws.Range("I28").Value = 1224.7826
$ws.Range("J28").Value = 364.66666
$ws.Range("K28").Value = 1224.7826
$ws.Range("L28").Value = 364.66666
$ws.Range("M28").Value = -739.7826
$ws.Range("N28").Value = -1334.66666

$ws.Range("H32").Value = 45459540
$ws.Range("I32").Value = 250002260
$ws.Range("J32").Value = 5605.6665
$ws.Range("K32").Value = 250002260
$ws.Range("L32").Value = 5605.6665
$ws.Range("M32").Value = -250001934
$ws.Range("N32").Value = -6257.6665

$ws.Range("H40").Value = 17640.834
$ws.Range("I40").Value = 15375.471
$ws.Range("K40").Value = 15375.471
$ws.Range("M40").Value = -15200.471

$ws.Range("H43").Value = 3198.3333
$ws.Range("I43").Value = 3197.5
$ws.Range("J43").Value = 3200
$ws.Range("K43").Value = 3197.5
$ws.Range("L43").Value = 3200
$ws.Range("M43").Value = -3128.5
$ws.Range("N43").Value = -3338

$ws.Range("H82").Value = 7727.143
$ws.Range("I82").Value = 2045
$ws.Range("K82").Value = 6135
$ws.Range("M82").Value = -5729

$ws.Range("H85").Value = 7727.143
$ws.Range("I85").Value = 2045
$ws.Range("K85").Value = 6135
$ws.Range("M85").Value = -4731

$ws.Range("H88").Value = 11249.75
$ws.Range("J88").Value = 12099.4
$ws.Range("L88").Value = 12099.4
$ws.Range("N88").Value = -12911.4

$ws.Range("H91").Value = 11249.75
$ws.Range("J91").Value = 12099.4
$ws.Range("L91").Value = 12099.4
$ws.Range("N91").Value = -14907.4

$ws.Range("H100").Value = 3083.2856
$ws.Range("I100").Value = 2358.875
$ws.Range("J100").Value = 4049.1667
$ws.Range("K100").Value = 2358.875
$ws.Range("L100").Value = 4049.1667
$ws.Range("M100").Value = -1817.875
$ws.Range("N100").Value = -5131.1667

$ws.Range("H112").Value = 3601.7632
$ws.Range("I112").Value = 1399.6666
$ws.Range("J112").Value = 3790.5144
$ws.Range("K112").Value = 4198.9998
$ws.Range("L112").Value = 11371.5432
$ws.Range("M112").Value = -3090.9998
$ws.Range("N112").Value = -13587.5432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13968.598
$ws.Range("I32").Value = 14499.982
$ws.Range("J32").Value = 11949.333
$ws.Range("K32").Value = 14499.982
$ws.Range("L32").Value = 11949.333
$ws.Range("M32").Value = -14212.982
$ws.Range("N32").Value = -12523.333

$ws.Range("H45").Value = 2755.4
$ws.Range("I45").Value = 1946.35
$ws.Range("J45").Value = 5991.6
$ws.Range("K45").Value = 1946.35
$ws.Range("L45").Value = 5991.6
$ws.Range("M45").Value = -1569.35
$ws.Range("N45").Value = -6745.6

$ws.Range("H97").Value = 820.84375
$ws.Range("I97").Value = 672.72
$ws.Range("K97").Value = 672.72
$ws.Range("M97").Value = -176.72

$ws.Range("H109").Value = 93999.664
$ws.Range("J109").Value = 93999.664
$ws.Range("L109").Value = 93999.664
$ws.Range("N109").Value = -96773.664

$ws.Range("H122").Value = 3173.0535
$ws.Range("I122").Value = 2157.0637
$ws.Range("K122").Value = 6471.1911
$ws.Range("M122").Value = -4021.1911

$ws.Range("H132").Value = 14488.652
$ws.Range("I132").Value = 16388.594
$ws.Range("K132").Value = 49165.78200000001
$ws.Range("M132").Value = -46635.78200000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1956.0714
$ws.Range("I99").Value = 1790.4166
$ws.Range("K99").Value = 1790.4166
$ws.Range("M99").Value = -292.4166

$ws.Range("H105").Value = 4204.1
$ws.Range("I105").Value = 3666.5
$ws.Range("J105").Value = 5010.5
$ws.Range("K105").Value = 3666.5
$ws.Range("L105").Value = 5010.5
$ws.Range("M105").Value = -1919.5
$ws.Range("N105").Value = -8504.5

$ws.Range("H107").Value = 2323.4187
$ws.Range("I107").Value = 1845
$ws.Range("J107").Value = 3559.3333
$ws.Range("K107").Value = 1845
$ws.Range("L107").Value = 3559.3333
$ws.Range("M107").Value = 75
$ws.Range("N107").Value = -7399.3333

$ws.Range("H134").Value = 2800.9333
$ws.Range("I134").Value = 2001.8572
$ws.Range("K134").Value = 6005.571599999999
$ws.Range("M134").Value = -3470.571599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 2000
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 2000
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 2000
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -2280

$ws.Range("H31").Value = 5250.4062
$ws.Range("I31").Value = 4200.5835
$ws.Range("K31").Value = 4200.5835
$ws.Range("M31").Value = -3905.5835

$ws.Range("H34").Value = 5250.4062
$ws.Range("I34").Value = 4200.5835
$ws.Range("K34").Value = 4200.5835
$ws.Range("M34").Value = -3998.5835

$ws.Range("H132").Value = 7102744.5
$ws.Range("I132").Value = 7943160.5
$ws.Range("K132").Value = 23829481.5
$ws.Range("M132").Value = -23826951.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1994.4375
$ws.Range("I114").Value = 625.6667
$ws.Range("K114").Value = 1877.0001
$ws.Range("M114").Value = 1376.9999

$ws.Range("H129").Value = 1047.1428
$ws.Range("I129").Value = 587.64703
$ws.Range("K129").Value = 1762.94109
$ws.Range("M129").Value = 3237.05891

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 803
$ws.Range("I107").Value = 564.3
$ws.Range("K107").Value = 564.3
$ws.Range("M107").Value = 1355.7

$ws.Range("H132").Value = 672258.1
$ws.Range("I132").Value = 138709.86
$ws.Range("K132").Value = 416129.58
$ws.Range("M132").Value = -413599.58

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7706
$ws.Range("I122").Value = 4443.1577
$ws.Range("J122").Value = 14594.223
$ws.Range("K122").Value = 13329.4731
$ws.Range("L122").Value = 43782.669
$ws.Range("M122").Value = -10879.4731
$ws.Range("N122").Value = -48682.669

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 112316.5
$ws.Range("J41").Value = 112316.5
$ws.Range("L41").Value = 112316.5
$ws.Range("N41").Value = -113096.5

$ws.Range("H55").Value = 12724.714
$ws.Range("I55").Value = 3773
$ws.Range("J55").Value = 16305.4
$ws.Range("K55").Value = 3773
$ws.Range("L55").Value = 16305.4
$ws.Range("M55").Value = -3496
$ws.Range("N55").Value = -16859.4

$ws.Range("H107").Value = 1365.28
$ws.Range("I107").Value = 1454.238
$ws.Range("K107").Value = 4362.714
$ws.Range("M107").Value = -2442.714

$ws.Range("H113").Value = 1656.25
$ws.Range("I113").Value = 1002.8889
$ws.Range("K113").Value = 3008.6667
$ws.Range("M113").Value = -838.6667000000002

$ws.Range("H122").Value = 6552.923
$ws.Range("I122").Value = 5031.6665
$ws.Range("J122").Value = 7856.857
$ws.Range("K122").Value = 15094.9995
$ws.Range("L122").Value = 23570.571
$ws.Range("M122").Value = -12644.9995
$ws.Range("N122").Value = -28470.571
